$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-14: refresh randomly generated credentials + candidate IDs (new iAuthor batch) ---
$ws.Range("A2").Value = 'TBZtO234'
$ws.Range("B2").Value = 231102272
$ws.Range("C2").Value = 'bujqqur42'
$ws.Range("D2").Value = 'D&gd4%8J'
$ws.Range("F2").Value = 'gxLCjAlv'
$ws.Range("G2").Value = 'Otxz'

$ws.Range("A3").Value = 'fqCMT400'
$ws.Range("B3").Value = 231102271
$ws.Range("C3").Value = 'uxufyye62'
$ws.Range("D3").Value = 'R%s2Eg5&'
$ws.Range("F3").Value = 'jvtzpDrF'
$ws.Range("G3").Value = 'cKAM'

$ws.Range("A4").Value = 'JxAty549'
$ws.Range("B4").Value = 231102270
$ws.Range("C4").Value = 'rdakwap42'
$ws.Range("D4").Value = 'h9RC%$g8'
$ws.Range("F4").Value = 'mcpFybdd'
$ws.Range("G4").Value = 'RPfG'

$ws.Range("A5").Value = 'EfTFE893'
$ws.Range("B5").Value = 231102269
$ws.Range("C5").Value = 'toqsdpq26'
$ws.Range("D5").Value = 'pD6g$V!9'
$ws.Range("F5").Value = 'uDLSYPUq'
$ws.Range("G5").Value = 'MjNp'

$ws.Range("A6").Value = 'aHmsF874'
$ws.Range("B6").Value = 231102268
$ws.Range("C6").Value = 'rkohqul24'
$ws.Range("D6").Value = 'U2xrG!6#'
$ws.Range("F6").Value = 'FSwIPkQh'
$ws.Range("G6").Value = 'kDLg'

$ws.Range("A7").Value = 'CztTs612'
$ws.Range("B7").Value = 231102267
$ws.Range("C7").Value = 'kzreddl96'
$ws.Range("D7").Value = 'ak25EB#!'
$ws.Range("F7").Value = 'KbzZSviF'
$ws.Range("G7").Value = 'XDIh'

$ws.Range("A8").Value = 'oGsfI795'
$ws.Range("B8").Value = 231102266
$ws.Range("C8").Value = 'fhberbb84'
$ws.Range("D8").Value = 'k$7!eVJ5'
$ws.Range("F8").Value = 'cgwMXmsX'
$ws.Range("G8").Value = 'qOuN'

$ws.Range("A9").Value = 'eMqwR237'
$ws.Range("B9").Value = 231102265
$ws.Range("C9").Value = 'xsmpscz70'
$ws.Range("D9").Value = 'j$3YS8p!'
$ws.Range("F9").Value = 'RcTOXtaD'
$ws.Range("G9").Value = 'EFmT'

$ws.Range("A10").Value = 'FCfdv958'
$ws.Range("B10").Value = 231102264
$ws.Range("C10").Value = 'eznsepg50'
$ws.Range("D10").Value = 'x3z4R!%P'
$ws.Range("F10").Value = 'JnxgkRYx'
$ws.Range("G10").Value = 'xCvU'

$ws.Range("A11").Value = 'oVkaK186'
$ws.Range("B11").Value = 231102263
$ws.Range("C11").Value = 'aqjofrp65'
$ws.Range("D11").Value = 'a5!4XPd&'
$ws.Range("F11").Value = 'vFbCZoyH'
$ws.Range("G11").Value = 'wnSD'

$ws.Range("A12").Value = 'MKzDW348'
$ws.Range("B12").Value = 231102262
$ws.Range("C12").Value = 'okrfvru60'
$ws.Range("D12").Value = 'd#4BM7x%'
$ws.Range("F12").Value = 'PrZYbQhZ'
$ws.Range("G12").Value = 'kVjH'

$ws.Range("A13").Value = 'CutFi497'
$ws.Range("B13").Value = 231102261
$ws.Range("C13").Value = 'tbedorh37'
$ws.Range("D13").Value = 'y85R!hV#'
$ws.Range("F13").Value = 'yTYpwhFp'
$ws.Range("G13").Value = 'tupz'

$ws.Range("A14").Value = 'zJnpd385'
$ws.Range("B14").Value = 231102260
$ws.Range("C14").Value = 'rrjyuyd18'
$ws.Range("D14").Value = 'qP#6%9Ub'
$ws.Range("F14").Value = 'hOEVXNad'
$ws.Range("G14").Value = 'XRFR'

# --- Append newly added row 15 (new iAuthor TC) ---
$ws.Range("A15").Value = 'WhGgh611'
$ws.Range("B15").Value = 231102259
$ws.Range("C15").Value = 'ywuojdr55'
$ws.Range("D15").Value = 'Xh8$U#6w'
$ws.Range("E15").Value = 'MR'
$ws.Range("F15").Value = 'LtOROZDr'
$ws.Range("G15").Value = 'BwVz'
$ws.Range("H15").Value = 'Candidate'

# Match the thin-border styling used by the other data rows (2-14)
$ws.Range("A15:H15").Borders.LineStyle = 1

# Grow the sheet dimension/selection to cover the newly added row
$null = $ws.Range("A1:H15").Select()
